$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 'BSV HalleAmmendorf'
$ws.Range("B46").Value = 7248441
$ws.Range("F46").Value = 'SV BlauWeiss Zorbau'
$ws.Range("G46").Value = 3
$ws.Range("H46").Value = 1
$ws.Range("I46").Value = 'H'
$ws.Range("J46").Value = 3.25
$ws.Range("K46").Value = 3.8
$ws.Range("L46").Value = 1.833
$ws.Range("M46").Value = 2.7
$ws.Range("N46").Value = 3.75
$ws.Range("O46").Value = 2.1
$ws.Range("Q46").Value = 1.85
$ws.Range("R46").Value = 1.95
$ws.Range("S46").Value = 3
$ws.Range("V46").Value = 1.7
$ws.Range("X46").Value = -1
$ws.Range("Y46").Value = 0.8500000000000001
$ws.Range("Z46").Value = -1
$ws.Range("B47").Value = 7248791
$ws.Range("F47").Value = 'FC Zuzenhausen'
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 3
$ws.Range("I47").Value = 'A'
$ws.Range("J47").Value = 2.75
$ws.Range("K47").Value = 4
$ws.Range("L47").Value = 2
$ws.Range("M47").Value = 2.75
$ws.Range("N47").Value = 4
$ws.Range("O47").Value = 2
$ws.Range("Q47").Value = 1.975
$ws.Range("R47").Value = 1.825
$ws.Range("S47").Value = 3.25
$ws.Range("V47").Value = -1
$ws.Range("X47").Value = 1
$ws.Range("Y47").Value = -1
$ws.Range("Z47").Value = 0.825
$ws.Range("B67").Value = 7423701
$ws.Range("F67").Value = 'VfB Sangerhausen'
$ws.Range("H67").Value = 2
$ws.Range("I67").Value = 'A'
$ws.Range("J67").Value = 2
$ws.Range("L67").Value = 2.55
$ws.Range("M67").Value = 2
$ws.Range("O67").Value = 2.6
$ws.Range("P67").Value = -0.25
$ws.Range("S67").Value = 3.25
$ws.Range("W67").Value = -1
$ws.Range("X67").Value = 1.6
$ws.Range("B68").Value = 7423702
$ws.Range("F68").Value = 'SV Fortuna Magdeburg'
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 'D'
$ws.Range("J68").Value = 1.8
$ws.Range("L68").Value = 3
$ws.Range("M68").Value = 1.8
$ws.Range("O68").Value = 3
$ws.Range("P68").Value = -0.5
$ws.Range("S68").Value = 3.5
$ws.Range("W68").Value = 3.5
$ws.Range("X68").Value = -1
$ws.Range("E95").Value = 'BSV HalleAmmendorf'
$ws.Range("E98").Value = 'SG RotWeiss Thalheim'
$ws.Range("F104").Value = 'Waldhof Mannheim II'
$ws.Range("E110").Value = 'SG RotWeiss Thalheim'
$ws.Range("E113").Value = 'SG RotWeiss Thalheim'
$ws.Range("F113").Value = 'BSV HalleAmmendorf'
$ws.Range("E119").Value = 'SG RotWeiss Thalheim'
$ws.Range("E123").Value = 'SG RotWeiss Thalheim'
$ws.Range("E132").Value = '1 FC Lok Stendal'
$ws.Range("E135").Value = 'SG RotWeiss Thalheim'
$ws.Range("E139").Value = 'BSV HalleAmmendorf'
$ws.Range("E146").Value = 'Waldhof Mannheim II'
